# Page 10 (Slide 10): fix the description of "$@" from
# "$@: Object Files" to "$@: the target file".
#
# The bullet's first run ("$@") keeps its bold/blue styling; the rest of
# the line ("<colon> Object Files") is rewritten in two steps so the
# resulting run layout matches the authored edit:
#   1) the tail words "Object Files" become "the target file"
#   2) the final word is re-applied so it lands in its own run
#      ("$@" / ": the target " / "file")

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(10)
$shp = $s.Shapes.Item(2)          # "Content Placeholder 2"

$tr    = $shp.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)     # "$@: Object Files"

# Replace ": Object Files" (everything after the "$@" run) in one go.
$rest = $para1.Characters(3, 14)
$rest.Text = ": the target file"

# Re-type just "file" so it ends up as its own trailing run, matching
# "$@" + ": the target " + "file".
$tail = $para1.Characters(16, 4)
$tail.Text = "file"
